$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily format the entire Price column as Text so numeric-looking
# strings (e.g. '32.20', '1.00') are preserved verbatim instead of being
# parsed/rounded as numbers by Excel's automatic type detection.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '34.609.98'
$ws.Range("E2").Value = '  +1.56%  '
$ws.Range("D3").Value = '1.793.06'
$ws.Range("E3").Value = '  +0.45%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = '223.02'
$ws.Range("E5").Value = '  -1.31%  '
$ws.Range("E6").Value = '  -1.23%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").Value = '32.20'
$ws.Range("E8").Value = '  +7.13%  '
$ws.Range("E9").Value = '  +0.28%  '
$ws.Range("D10").Value = '0.0687'
$ws.Range("E10").Value = '  +2.73%  '
$ws.Range("D11").Value = '0.0936'
$ws.Range("E11").Value = '  +1.49%  '
$ws.Range("D12").Value = '2.051.67'
$ws.Range("E12").Value = '  +0.63%  '
$ws.Range("D13").Value = '1.775.97'
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("E14").Value = '  +7.43%  '
$ws.Range("D15").Value = '34.637.77'
$ws.Range("E15").Value = '  +1.84%  '
$ws.Range("D16").Value = '0.630'
$ws.Range("E16").Value = '  +0.90%  '
$ws.Range("E17").Value = '  +2.18%  '
$ws.Range("D18").Value = '68.52'
$ws.Range("E18").Value = '  -0.32%  '
$ws.Range("D19").Value = '254.16'
$ws.Range("E19").Value = '  +0.84%  '
$ws.Range("D20").Value = '0.0₃0789'
$ws.Range("E20").Value = '  +6.10%  '
$ws.Range("D21").Value = '1.00'
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("D22").Value = '10.43'
$ws.Range("E22").Value = '  +0.54%  '
$ws.Range("D23").Value = '4.16'
$ws.Range("E23").Value = '  -0.73%  '
$ws.Range("E24").Value = '  +0.51%  '
$ws.Range("D25").Value = '159.86'
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("E26").Value = '  -0.91%  '
$ws.Range("E27").Value = '  +1.21%  '
$ws.Range("E28").Value = '  -0.47%  '
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = '0.0516'
$ws.Range("E30").Value = '  +0.23%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '3.75'
$ws.Range("E31").Value = '  -2.36%  '
$ws.Range("E32").Value = '  -0.27%  '
$ws.Range("E33").Value = '  -0.85%  '
$ws.Range("E34").Value = '  +1.29%  '
$ws.Range("D35").Value = '1.435.97'
$ws.Range("E35").Value = '  -4.39%  '
$ws.Range("E36").Value = '  -1.07%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '0.636'
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.0191'
$ws.Range("E38").Value = '  +2.55%  '
$ws.Range("D39").Value = '84.30'
$ws.Range("E39").Value = '  +0.94%  '
$ws.Range("E40").Value = '  +3.75%  '
$ws.Range("E41").Value = '  +0.21%  '
$ws.Range("D42").Value = '0.910'
$ws.Range("E42").Value = '  +1.44%  '
$ws.Range("E43").Value = '  +1.31%  '
$ws.Range("D44").Value = '6.04'
$ws.Range("E44").Value = '  +5.44%  '
$ws.Range("E45").Value = '  -1.40%  '
$ws.Range("D46").Value = '0.0494'
$ws.Range("E46").Value = '  -3.76%  '
$ws.Range("D47").Value = '1.949.40'
$ws.Range("E47").Value = '  +0.54%  '
$ws.Range("E48").Value = '  +1.52%  '
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("D50").Value = '103.77'
$ws.Range("E50").Value = '  +5.94%  '
$ws.Range("D51").Value = '49.81'
$ws.Range("E51").Value = '  -3.42%  '

# Restore the default (General) style on the whole Price column so no
# stray text-format style is left attached to any of these cells.
$ws.Range("D2:D51").Style = "Normal"

Write-Output "applied cryptos update"